$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.737.76"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").Value = "'1.680.92"
$ws.Range("E3").Value = "  -1.52%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").Value = "'314.94"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("E6").Value = "  +0.38%  "

$ws.Range("D7").Value = "'0.3925"
$ws.Range("E7").Value = "  -1.66%  "

$ws.Range("D8").Value = "'0.3962"
$ws.Range("E8").Value = "  -2.74%  "

$ws.Range("D9").Value = "'1.002"
$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("D10").Value = "'51.76"
$ws.Range("E10").Value = "  -3.79%  "

$ws.Range("D11").Value = "'1.400"
$ws.Range("E11").Value = "  -5.87%  "

$ws.Range("D12").Value = "'0.08639"
$ws.Range("E12").Value = "  -2.13%  "

$ws.Range("D13").Value = "'25.30"
$ws.Range("E13").Value = "  -3.84%  "

$ws.Range("D14").Value = "'7.335"
$ws.Range("E14").Value = "  -2.35%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.758"
$ws.Range("E15").Value = "  -5.02%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.00001313"
$ws.Range("E16").Value = "  -3.64%  "

$ws.Range("D17").Value = "'1.692.31"
$ws.Range("E17").Value = "  -1.32%  "

$ws.Range("D18").Value = "'93.17"
$ws.Range("E18").Value = "  -3.96%  "

$ws.Range("D19").Value = "'0.07095"
$ws.Range("E19").Value = "  -1.17%  "

$ws.Range("D20").Value = "'20.23"
$ws.Range("E20").Value = "  -4.59%  "

$ws.Range("D21").Value = "'7.052"
$ws.Range("E21").Value = "  -3.12%  "

$ws.Range("E22").Value = "  +0.50%  "

$ws.Range("D23").Value = "'13.92"
$ws.Range("E23").Value = "  -3.52%  "

$ws.Range("D24").Value = "'24.729.04"
$ws.Range("E24").Value = "  -0.59%  "

$ws.Range("D25").Value = "'2.349"
$ws.Range("E25").Value = "  +0.99%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.828"
$ws.Range("E26").Value = "  -3.87%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'23.35"
$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").Value = "'162.79"
$ws.Range("E28").Value = "  -2.28%  "

$ws.Range("D29").Value = "'5.829"
$ws.Range("E29").Value = "  -6.28%  "

$ws.Range("D30").Value = "'146.72"
$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("D31").Value = "'7.883"
$ws.Range("E31").Value = "  -6.41%  "

$ws.Range("D32").Value = "'2.377"
$ws.Range("E32").Value = "  +6.27%  "

$ws.Range("D33").Value = "'1.877.67"
$ws.Range("E33").Value = "  -1.14%  "

$ws.Range("D34").Value = "'0.08392"
$ws.Range("E34").Value = "  -4.64%  "

$ws.Range("D35").Value = "'0.03056"
$ws.Range("E35").Value = "  -5.14%  "

$ws.Range("D36").Value = "'6.953"
$ws.Range("E36").Value = "  -3.95%  "

$ws.Range("D37").Value = "'1.000"
$ws.Range("E37").Value = "  -3.18%  "

$ws.Range("D38").Value = "'0.2791"
$ws.Range("E38").Value = "  -3.36%  "

$ws.Range("D39").Value = "'0.09442"
$ws.Range("E39").Value = "  +1.15%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'10.62"
$ws.Range("E40").Value = "  -2.62%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.543"
$ws.Range("E41").Value = "  +4.99%  "

$ws.Range("D42").Value = "'0.7926"
$ws.Range("E42").Value = "  -8.06%  "

$ws.Range("D43").Value = "'13.49"
$ws.Range("E43").Value = "  -5.38%  "

$ws.Range("D44").Value = "'16.58"
$ws.Range("E44").Value = "  -5.99%  "

$ws.Range("D45").Value = "'0.7139"
$ws.Range("E45").Value = "  -4.66%  "

$ws.Range("D46").Value = "'2.563"
$ws.Range("E46").Value = "  -5.52%  "

$ws.Range("D47").Value = "'4.189"
$ws.Range("E47").Value = "  -1.35%  "

$ws.Range("D48").Value = "'0.08657"
$ws.Range("E48").Value = "  +3.49%  "

$ws.Range("E49").Value = "  +0.42%  "

$ws.Range("D50").Value = "'1.341"
$ws.Range("E50").Value = "  -4.37%  "

$ws.Range("D51").Value = "'137.83"
$ws.Range("E51").Value = "  -2.36%  "

